# Re-orders the weekly price records (columns D, L, M, N, O, P, Q, S, T)
# across rows 2-9 while keeping the market/region/product metadata columns
# (A, B, C, E-K, R) untouched, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for the columns that move, keyed by row.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")
$snapshot = @{}
foreach ($r in 2..9) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

# Destination row -> source row mapping (after <- before).
$map = @{
    2 = 8
    3 = 7
    4 = 6
    5 = 9
    6 = 4
    7 = 5
    8 = 2
    9 = 3
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
